# Add "Spain" and "Italy" test-data sheets, cloned from the existing
# "Norway" sheet template, then reorder so Italy precedes Spain and
# Italy becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook
$norway = $wb.Worksheets.Item("Norway")

# 1) Clone Norway -> Spain (placed right after Norway) and fill its data.
$norway.Copy($null, $norway)
$spain = $wb.Worksheets.Item(6)
$spain.Name = "Spain"
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3442/T2125"

# 2) Clone Norway -> Italy (placed right after Spain) and fill its data.
$norway.Copy($null, $spain)
$italy = $wb.Worksheets.Item(7)
$italy.Name = "Italy"
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3443/T1972"

# 3) Reorder: Italy should come before Spain.
$italy.Move($spain, $null)

# Re-fetch handles by name since Move() can invalidate index-bound
# references to the sheets that shifted position.
$spain = $wb.Worksheets.Item("Spain")
$italy = $wb.Worksheets.Item("Italy")

# 4) Restore each sheet's last selection and make Italy the active tab.
$spain.Range("A3").Select()
$italy.Range("A10").Select()
$italy.Activate()
